$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "level_id" header in A1 was renamed to "user_id".
$ws.Range("A1").Value = "user_id"

# Selection had been left on E2 from a prior session; reset it back to A1
# (the sheet's natural/default selection) to match the saved view state.
$ws.Range("A1").Select() | Out-Null
